$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at CJ (column 88). This shifts the existing
# "nom" (CJ -> CK) and "url_produit" (CK -> CL) columns one to the right,
# and grows the used range from A1:CK206 to A1:CL206.
$ws.Columns("CJ").Insert()

# New CJ1 header is a fresh snapshot timestamp.
$ws.Range("CJ1").Value = "2026-01-31 16:13:48"

# Every data row's new CJ cell repeats the last known price (the value
# that was in CI before the insert, which is still in CI after the
# insert since the insert happened at CJ). Rows with no tracked price
# (CI blank) are left blank too.
for ($r = 2; $r -le 206; $r++) {
    $lastPrice = $ws.Cells.Item($r, 87).Value()
    if ($lastPrice -ne "") {
        $ws.Cells.Item($r, 88).Value = $lastPrice
    }
}
